# Vendor Pincode Mapping Table Cleaning
#
# The template originally had 10 columns per row:
#   A: vendor_name   F: area
#   B: vendor_id     G: pincode
#   C: appliance     H: region
#   D: appliance_id  I: city
#   E: brand         J: state
# (row 1 holds the lower-case placeholder names, row 2 holds the
# "{order:...}" merge tags)
#
# The cleaned template only keeps: vendor_name, appliance, pincode, city,
# state -- i.e. original columns B, D, E, F, H must be removed and the
# remaining columns shifted left so they end up as A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToLeft = -4159

# Delete the unwanted columns from right to left (within just the used
# rows, 1:2) so remaining column letters of not-yet-deleted columns stay
# valid while iterating.
$ws.Range("H1:H2").Delete($xlShiftToLeft) | Out-Null   # region
$ws.Range("F1:F2").Delete($xlShiftToLeft) | Out-Null   # area
$ws.Range("E1:E2").Delete($xlShiftToLeft) | Out-Null   # brand
$ws.Range("D1:D2").Delete($xlShiftToLeft) | Out-Null   # appliance_id
$ws.Range("B1:B2").Delete($xlShiftToLeft) | Out-Null   # vendor_id

# Leave the cursor where the author left it after cleaning up the sheet.
$ws.Range("H8").Select() | Out-Null
